# Weekly fruit/vegetable price update: insert a new weekly record row
# right above the previous "Provincia de Limarí" / Haba entry (row 152),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 152 (shifts rows 152..168 -> 153..169).
$ws.Rows.Item(152).Insert()

# Populate the new row with the latest weekly data point.
$ws.Range('A152').Value = 4
$ws.Range('B152').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C152').Value = 'Los Lagos'
$ws.Range('D152').Value = 45212
$ws.Range('E152').Value = 10
$ws.Range('F152').Value = 100112026
$ws.Range('G152').Value = 'Haba'
$ws.Range('H152').Value = 'Sin especificar'
$ws.Range('I152').Value = 'Primera'
$ws.Range('J152').Value = 120
$ws.Range('K152').Value = 16000
$ws.Range('L152').Value = 16000
$ws.Range('M152').Value = 16000
$ws.Range('N152').Value = '$/saco 25 kilos'
$ws.Range('O152').Value = 'Provincia de Limarí'
$ws.Range('P152').Value = 640
$ws.Range('Q152').Value = 25
$ws.Range('R152').Value = 'Hortaliza'
